$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Résultats" column in E, preserving the original optimal
# solution values there, while D now holds a new (integer) solution.

# Copy style of D1 to E1, D2 to E2, D3:D6 to E3:E6 so the new column
# matches the look of the existing "Valeurs" column.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("D3:D6").Copy()
$ws.Range("E3:E6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header for new column E
$ws.Range("E1").Value = "Résultats"

# Move the original solution values from D3/D4 into E3/E4
$ws.Range("E3").Value = 1.4545454545454546
$ws.Range("E4").Value = 1.9090909090909092
$ws.Range("E6").Formula = "=4*E3+3*E4"

# Set new solution values in D3/D4
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 1

# Update constraint formulas to reference the original (now in E) solution
$ws.Range("D9").Formula = "=3*E3+4*E4"
$ws.Range("D10").Formula = "=7*E3+2*E4"

# Update selection / active cell
$ws.Range("F5").Select()

# Update solver defined names to point at the new E column cells
$wb.Names.Item("solver_adj").RefersTo = "=Feuil1!`$E`$3,Feuil1!`$E`$4"
$wb.Names.Item("solver_opt").RefersTo = "=Feuil1!`$E`$6"
